# Change the table style applied to the table on slide 16 from
# "{76841C99-F89D-4C46-A836-46F1A0A4379D}" to
# "{1953C583-03F8-4023-8846-661647F27689}" (a:tblPr/a:tableStyleId).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTable) {
        $sh.Table.ApplyStyle("{1953C583-03F8-4023-8846-661647F27689}")
    }
}
